$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at A, shifting existing data (and the formula's
# relative references) one column to the right.
$ws.Columns("A").Insert()

# Row 2 tweaks (counts).
$ws.Range("C2").Value = 1
$ws.Range("E2").Value = 2

# Row 3 content updates (XP->X, INC->ic, Usuario->usua).
$ws.Range("C3").Value = "X"
$ws.Range("E3").Value = "ic"
$ws.Range("I3").Value = "usua"

# New "objeto" header column (A) with Tabla / Store rows.
$ws.Range("A1").Value = "objeto"
$ws.Range("A3").Value = "Tabla"
$ws.Range("A4").Value = "Store"

# Row 4 content updates.
$ws.Range("B4").Value = "SP"
$ws.Range("C4").Value = "X"
$ws.Range("E4").Value = "ic"
$ws.Range("H4").ClearContents()
$ws.Range("I4").Value = "Usuario"

# Rebuild the CONCATENATE formulas over the shifted columns (B..I) into J.
$ws.Range("J3").Formula = "=CONCATENATE(B3,C3,D3,E3,F3,G3,H3,I3)"
$ws.Range("J4").Formula = "=CONCATENATE(B4,C4,D4,E4,F4,G4,H4,I4)"

$ws.Range("B4").Select()
